$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (col D) and Volume(1h) (col E) values
# Cells whose new Price text would otherwise be auto-parsed as a number
# are forced to the Text number format first, so the literal string is kept.

$ws.Cells.Item(2, 4).Value = '29.940.60'
$ws.Cells.Item(2, 5).Value = '  +0.37%  '

$ws.Cells.Item(3, 4).Value = '1.892.80'
$ws.Cells.Item(3, 5).Value = '  -0.32%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.000'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.8311'
$ws.Cells.Item(5, 5).Value = '  +8.41%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '241.76'
$ws.Cells.Item(6, 5).Value = '  +0.61%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.000'
$ws.Cells.Item(7, 5).Value = '  -0.03%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3249'
$ws.Cells.Item(8, 5).Value = '  +6.58%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '26.66'
$ws.Cells.Item(9, 5).Value = '  +5.25%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.07027'
$ws.Cells.Item(10, 5).Value = '  +2.79%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.08027'
$ws.Cells.Item(11, 5).Value = '  +0.60%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.7478'
$ws.Cells.Item(12, 5).Value = '  +1.35%  '

$ws.Cells.Item(13, 4).Value = '1.895.74'
$ws.Cells.Item(13, 5).Value = '  +0.05%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.200'
$ws.Cells.Item(14, 5).Value = '  +0.63%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '92.30'
$ws.Cells.Item(15, 5).Value = '  +1.36%  '

$ws.Cells.Item(16, 4).Value = '29.937.90'
$ws.Cells.Item(16, 5).Value = '  +0.34%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '14.04'
$ws.Cells.Item(17, 5).Value = '  +1.92%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '5.921'
$ws.Cells.Item(18, 5).Value = '  +0.43%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '243.55'
$ws.Cells.Item(19, 5).Value = '  -0.54%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.000007758'
$ws.Cells.Item(20, 5).Value = '  +0.78%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.000'
$ws.Cells.Item(21, 5).Value = '  +0.04%  '

$ws.Cells.Item(22, 4).Value = '2.152.03'
$ws.Cells.Item(22, 5).Value = '  +1.03%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.001'
$ws.Cells.Item(23, 5).Value = '  -0.05%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.918'
$ws.Cells.Item(24, 5).Value = '  +0.06%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.1596'
$ws.Cells.Item(25, 5).Value = '  +24.10%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '167.73'
$ws.Cells.Item(26, 5).Value = '  +0.50%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.192'
$ws.Cells.Item(27, 5).Value = '  -0.67%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.84'
$ws.Cells.Item(28, 5).Value = '  +0.97%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.081'
$ws.Cells.Item(29, 5).Value = '  +2.42%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.371'
$ws.Cells.Item(30, 5).Value = '  -1.86%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.513'
$ws.Cells.Item(31, 5).Value = '  -0.07%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.260'
$ws.Cells.Item(32, 5).Value = '  -0.17%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.05628'
$ws.Cells.Item(33, 5).Value = '  +6.84%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.069'
$ws.Cells.Item(34, 5).Value = '  +0.01%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.276'
$ws.Cells.Item(35, 5).Value = '  +2.30%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.7324'
$ws.Cells.Item(36, 5).Value = '  +1.01%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.720'
$ws.Cells.Item(37, 5).Value = '  +0.03%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01910'
$ws.Cells.Item(38, 5).Value = '  +0.04%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.778'
$ws.Cells.Item(39, 5).Value = '  +0.00%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.4419'
$ws.Cells.Item(40, 5).Value = '  +0.44%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '71.90'
$ws.Cells.Item(41, 5).Value = '  -0.11%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.934'
$ws.Cells.Item(42, 5).Value = '  -4.13%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.8429'
$ws.Cells.Item(43, 5).Value = '  +1.15%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.000'
$ws.Cells.Item(44, 5).Value = '  -0.06%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.888'
$ws.Cells.Item(45, 5).Value = '  +0.48%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '101.32'
$ws.Cells.Item(46, 5).Value = '  +1.65%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.596'
$ws.Cells.Item(47, 5).Value = '  -0.02%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '9.699'
$ws.Cells.Item(48, 5).Value = '  -0.25%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '990.17'
$ws.Cells.Item(49, 5).Value = '  +9.06%  '

$ws.Cells.Item(50, 4).Value = '2.048.60'
$ws.Cells.Item(50, 5).Value = '  +0.63%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '36.01'
$ws.Cells.Item(51, 5).Value = '  -0.19%  '
